$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$row = $tbl.Rows.Item(12)
$row.Cells.Item(1).Range.Text = "07.10.2022"
$row.Cells.Item(2).Range.Text = "0,5"
$row.Cells.Item(3).Range.Text = "Koodin muokkaus, luokkadiagrammin päivitys"
